$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.660.98"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.75%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.923.86"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.39%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.012"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.44"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4817"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4054"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08099"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.65%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.002"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "23.48"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.969.14"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.991"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.192"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.11"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06849"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.012"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.45%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001028"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.54"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "29.674.36"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.552"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.77"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("E24").Value = "  -0.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.200.51"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.578"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.82"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +0.62%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.85"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -0.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.068"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "120.30"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.001"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.87%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09591"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.530"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.399"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.74%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.543"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06520"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +6.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02258"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.197"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5880"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "10.65"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.864"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1831"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.11%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.468"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.278"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.27"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -1.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.07463"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5509"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.962"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "116.33"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.400"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -1.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "71.93"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.72%  "
